# Apply the "reverse direct relation in references" fix to the
# hello_world_pump.xlsx rules workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Properties sheet: fix Container references to use the `cdf_cdm:`
# space prefix (they previously referenced bare, unprefixed container
# names), and fix the `livesIn` property to point at the `Facility`
# view instead of `CogniteAsset`. Also give `weight` its unit and drop
# the now-unused `weightUnit` row.
# ---------------------------------------------------------------------
$props = $wb.Worksheets.Item("Properties")

# Documentation.name -> container CogniteDescribable
$props.Range("K3").Value = "cdf_cdm:CogniteDescribable"

# Documentation.relatedPumps -> container CogniteFile
$props.Range("K4").Value = "cdf_cdm:CogniteFile"

# Facility.desc -> container CogniteDescribable
$props.Range("K7").Value = "cdf_cdm:CogniteDescribable"

# Facility.name -> container CogniteDescribable
$props.Range("K9").Value = "cdf_cdm:CogniteDescribable"

# Pump.livesIn -> value type Facility (was CogniteAsset), container cdf_cdm:CogniteAsset
$props.Range("F13").Value = "Facility"
$props.Range("K13").Value = "cdf_cdm:CogniteAsset"

# Pump.name -> container CogniteDescribable
$props.Range("K14").Value = "cdf_cdm:CogniteDescribable"

# Pump.weight -> value type carries the unit now
$props.Range("F17").Value = "float64(unit=mass:kilogm)"

# Pump.weightUnit row is no longer needed - remove it entirely (the
# row below, Pump.year, shifts up to take its place).
$props.Rows.Item(18).Delete()

# ---------------------------------------------------------------------
# Views sheet: add the CogniteTimeSeries view definition used by the
# Pump.pressure / Pump.temperature properties.
# ---------------------------------------------------------------------
$views = $wb.Worksheets.Item("Views")

$views.Range("A7").Value = "cdf_cdm:CogniteTimeSeries(version=v1)"
$views.Range("B7").Value = "Time series"
$views.Range("C7").Value = "Represents a series of data points in time order."
$views.Range("D7").Value = "cdf_cdm:CogniteDescribable(version=v1),cdf_cdm:CogniteSourceable(version=v1)"
$views.Range("F7").Value = $true
